$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2731.682
$ws.Range("I40").Value = 2343.889
$ws.Range("J40").Value = 3000.1538
$ws.Range("K40").Value = 2343.889
$ws.Range("L40").Value = 3000.1538
$ws.Range("M40").Value = -2168.889
$ws.Range("N40").Value = -3350.1538
$ws.Range("H58").Value = 4298.5
$ws.Range("I58").Value = 372
$ws.Range("J58").Value = 8225
$ws.Range("K58").Value = 1116
$ws.Range("L58").Value = 24675
$ws.Range("M58").Value = -966
$ws.Range("N58").Value = -24975
$ws.Range("H86").Value = 4006093
$ws.Range("I86").Value = 4024.1714
$ws.Range("J86").Value = 13344253
$ws.Range("K86").Value = 4024.1714
$ws.Range("L86").Value = 13344253
$ws.Range("M86").Value = -2901.1714
$ws.Range("N86").Value = -13346499
$ws.Range("H89").Value = 4006093
$ws.Range("I89").Value = 4024.1714
$ws.Range("J89").Value = 13344253
$ws.Range("K89").Value = 20120.857
$ws.Range("L89").Value = 66721265
$ws.Range("M89").Value = -14504.857
$ws.Range("N89").Value = -66732497
$ws.Range("H113").Value = 125004400
$ws.Range("I113").Value = 500003000
$ws.Range("J113").Value = 4866.8335
$ws.Range("K113").Value = 500003000
$ws.Range("L113").Value = 4866.8335
$ws.Range("M113").Value = -499999746
$ws.Range("N113").Value = -11374.8335
$ws.Range("H132").Value = 3445.4048
$ws.Range("I132").Value = 3107.7896
$ws.Range("K132").Value = 9323.3688
$ws.Range("M132").Value = -6793.3688
$ws.Range("H135").Value = 71430264
$ws.Range("I135").Value = 76924850
$ws.Range("J135").Value = 700
$ws.Range("K135").Value = 692323650
$ws.Range("L135").Value = 6300
$ws.Range("M135").Value = -692321115
$ws.Range("N135").Value = -11370
$ws.Range("H141").Value = 13379.5
$ws.Range("I141").Value = 1499
$ws.Range("K141").Value = 4497
$ws.Range("M141").Value = 683

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1812.1818
$ws.Range("I2").Value = 1748.2222
$ws.Range("K2").Value = 1748.2222
$ws.Range("M2").Value = -1635.2222
$ws.Range("H45").Value = 2737.4
$ws.Range("I45").Value = 2499.4
$ws.Range("K45").Value = 2499.4
$ws.Range("M45").Value = -2122.4
$ws.Range("H74").Value = 2870.7334
$ws.Range("I74").Value = 2947
$ws.Range("J74").Value = 2718.2
$ws.Range("K74").Value = 2947
$ws.Range("L74").Value = 2718.2
$ws.Range("M74").Value = -2073
$ws.Range("N74").Value = -4466.2
$ws.Range("H77").Value = 2870.7334
$ws.Range("I77").Value = 2947
$ws.Range("J77").Value = 2718.2
$ws.Range("K77").Value = 14735
$ws.Range("L77").Value = 13591
$ws.Range("M77").Value = -10367
$ws.Range("N77").Value = -22327
$ws.Range("H97").Value = 257
$ws.Range("I97").Value = 253.78947
$ws.Range("J97").Value = 287.5
$ws.Range("K97").Value = 253.78947
$ws.Range("L97").Value = 287.5
$ws.Range("M97").Value = 242.21053
$ws.Range("N97").Value = -1279.5
$ws.Range("H116").Value = 1812.1818
$ws.Range("I116").Value = 1748.2222
$ws.Range("K116").Value = 1748.2222
$ws.Range("M116").Value = 545.7778000000001
$ws.Range("H122").Value = 11113957
$ws.Range("I122").Value = 15153996
$ws.Range("K122").Value = 45461988
$ws.Range("M122").Value = -45459538
$ws.Range("H132").Value = 37039340
$ws.Range("I132").Value = 40002350
$ws.Range("J132").Value = 1724
$ws.Range("K132").Value = 120007050
$ws.Range("L132").Value = 5172
$ws.Range("M132").Value = -120004520
$ws.Range("N132").Value = -10232

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1812.1818
$ws.Range("I3").Value = 1748.2222
$ws.Range("K3").Value = 1748.2222
$ws.Range("M3").Value = -1634.2222
$ws.Range("H22").Value = 213.25
$ws.Range("I22").Value = 213.25
$ws.Range("K22").Value = 213.25
$ws.Range("M22").Value = -40.25
$ws.Range("H94").Value = 3527
$ws.Range("J94").Value = 1713.1818
$ws.Range("L94").Value = 1713.1818
$ws.Range("N94").Value = -2615.1818

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 999.5
$ws.Range("I16").Value = 999.5
$ws.Range("K16").Value = 999.5
$ws.Range("M16").Value = -712.5
$ws.Range("H31").Value = 2618.0127
$ws.Range("I31").Value = 1242.8064
$ws.Range("K31").Value = 1242.8064
$ws.Range("M31").Value = -947.8063999999999
$ws.Range("H34").Value = 2618.0127
$ws.Range("I34").Value = 1242.8064
$ws.Range("K34").Value = 1242.8064
$ws.Range("M34").Value = -1040.8064
$ws.Range("H107").Value = 1970.7333
$ws.Range("I107").Value = 2022.2142
$ws.Range("K107").Value = 2022.2142
$ws.Range("M107").Value = -102.2141999999999
$ws.Range("H113").Value = 999.5
$ws.Range("I113").Value = 999.5
$ws.Range("K113").Value = 999.5
$ws.Range("M113").Value = 1170.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 83
$ws.Range("I23").Value = 75
$ws.Range("J23").Value = 87
$ws.Range("K23").Value = 225
$ws.Range("L23").Value = 261
$ws.Range("M23").Value = 10
$ws.Range("N23").Value = -731

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 19166.666
$ws.Range("J20").Value = 19166.666
$ws.Range("L20").Value = 19166.666
$ws.Range("N20").Value = -19656.666
$ws.Range("H24").Value = 50000
$ws.Range("J24").Value = 50000
$ws.Range("L24").Value = 50000
$ws.Range("N24").Value = -50346
$ws.Range("H97").Value = 2864.75
$ws.Range("I97").Value = 2819.6667
$ws.Range("K97").Value = 2819.6667
$ws.Range("M97").Value = -2323.6667
$ws.Range("H102").Value = 2180.1614
$ws.Range("I102").Value = 1109
$ws.Range("K102").Value = 1109
$ws.Range("M102").Value = 513
$ws.Range("H113").Value = 2831.2104
$ws.Range("I113").Value = 1362
$ws.Range("K113").Value = 1362
$ws.Range("M113").Value = 808
$ws.Range("H126").Value = 8542.474
$ws.Range("I126").Value = 11490.272
$ws.Range("K126").Value = 34470.81600000001
$ws.Range("M126").Value = -32000.81600000001

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15154723
$ws.Range("I7").Value = 20836002
$ws.Range("K7").Value = 20836002
$ws.Range("M7").Value = -20835890
$ws.Range("H21").Value = 5671.3335
$ws.Range("J21").Value = 5671.3335
$ws.Range("L21").Value = 5671.3335
$ws.Range("N21").Value = -6019.3335
$ws.Range("H126").Value = 15154723
$ws.Range("I126").Value = 20836002
$ws.Range("K126").Value = 62508006
$ws.Range("M126").Value = -62505536
$ws.Range("H132").Value = 5556.727
$ws.Range("I132").Value = 3515.111
$ws.Range("J132").Value = 8006.6665
$ws.Range("K132").Value = 10545.333
$ws.Range("L132").Value = 24019.9995
$ws.Range("M132").Value = -8015.332999999999
$ws.Range("N132").Value = -29079.9995

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9523.809999999999
$ws.Range("H64").Value = 26700
$ws.Range("I64").Value = 26700
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 26700
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -26452
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 26700
$ws.Range("I67").Value = 26700
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 26700
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -25842
$ws.Range("N67").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 6097.2104
$ws.Range("I132").Value = 6143.3076
$ws.Range("K132").Value = 18429.9228
$ws.Range("M132").Value = -15899.9228
$ws.Range("H136").Value = 5333.615
$ws.Range("I136").Value = 2224
$ws.Range("K136").Value = 6672
$ws.Range("M136").Value = -4122
